$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.064.31"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.834.10"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D4").Value = "'0.9991"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'242.74"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").Value = "'0.6174"
$ws.Range("E6").Value = "  -2.45%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.07452"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "'0.2927"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "'23.02"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "1.838.12"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "'4.996"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").Value = "'0.6722"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "'82.72"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").Value = "'0.000009170"
$ws.Range("E16").Value = "  -4.57%  "
$ws.Range("D17").Value = "'5.903"
$ws.Range("E17").Value = "  -2.66%  "
$ws.Range("D18").Value = "29.028.65"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "2.078.01"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "'239.53"
$ws.Range("E20").Value = "  +5.69%  "
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'7.214"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'159.35"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").Value = "'0.1414"
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "'17.84"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").Value = "'1.500"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "'4.143"
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.05570"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.105"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").Value = "'1.207"
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.841"
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7421"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").Value = "'2.772"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").Value = "'0.01783"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "1.212.73"
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("D41").Value = "'6.428"
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("D42").Value = "'0.8989"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").Value = "'0.9999"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'101.48"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "1.977.00"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'65.55"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "'0.5086"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("E48").Value = "  -2.56%  "
$ws.Range("D49").Value = "'0.4067"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "'9.094"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("E51").Value = "  +0.54%  "
